# feat: add 2022-Q3 data
#
# 1) Insert a new "2022-Q3" quarter sheet (cloned from the "2022-Q2" sheet so it
#    inherits the same sheet/structure boilerplate), positioned right after the
#    "总计" summary sheet and before "2022-Q2".
# 2) Fill it with the new fund-holding rows for 2022-Q3.
# 3) Insert a matching row at the top of the "总计" summary sheet and bump the
#    running index column for the rows that shifted down.

function Set-TextCell($ws, $row, $col, $text) {
    # Force literal-text storage (no numeric auto-coercion of things like
    # "167506" / "0.82" / "93.20") while keeping the cell's style at the
    # workbook default (no stray NumberFormat residue).
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: clone "2022-Q2" (currently sheet index 2) to create "2022-Q3" before it
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Drop the two data rows the clone inherited from 2022-Q2 (rows 2..5); we only
# need two data rows (2..3) for 2022-Q3, which we overwrite below.
$q3.Rows.Item(4).Delete()
$q3.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# Step 2: populate the 2022-Q3 sheet
# ---------------------------------------------------------------------------
# Header cells already carry the bold/border header style inherited from the
# cloned sheet (s="2") - plain .Value assignment keeps it (headers are not
# numeric-looking, so there is no auto-coercion risk here).
$q3.Cells.Item(1,2).Value = "基金代码"
$q3.Cells.Item(1,3).Value = "基金名称"
$q3.Cells.Item(1,4).Value = "基金规模"
$q3.Cells.Item(1,5).Value = "股票总仓位"
$q3.Cells.Item(1,6).Value = "仓位占比"
$q3.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value = "仓位排名"

$q3.Cells.Item(2,1).Value = 0
Set-TextCell $q3 2 2 "167506"
Set-TextCell $q3 2 3 "安信深圳科技指数（LOF）A"
Set-TextCell $q3 2 4 "0.82"
Set-TextCell $q3 2 5 "93.20"
Set-TextCell $q3 2 6 "2.82"
Set-TextCell $q3 2 7 "0.0231"
$q3.Cells.Item(2,8).Value = 10

$q3.Cells.Item(3,1).Value = 1
Set-TextCell $q3 3 2 "167507"
Set-TextCell $q3 3 3 "安信深圳科技指数（LOF）C"
Set-TextCell $q3 3 4 "0.30"
Set-TextCell $q3 3 5 "93.20"
Set-TextCell $q3 3 6 "2.82"
Set-TextCell $q3 3 7 "0.0085"
$q3.Cells.Item(3,8).Value = 10

# ---------------------------------------------------------------------------
# Step 3: add the 2022-Q3 row to the "总计" summary sheet
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$total.Rows.Item(2).Insert()

# New row's B:D cells come back with a stray inherited style - match the
# original sheet, where only column A carries the bold/border style.
$total.Range("B2:D2").ClearFormats()

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122)

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 2
$total.Cells.Item(2,4).Value = 0.03

# Bump the 0-based running index for every row that shifted down one slot.
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(7,1).Value = 5
